$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

# --- Swap the "description" (B) and "type" (E) columns, including the header row ---
for ($r = 1; $r -le 16; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 2).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $bVal
}

# --- Swap the column widths that belonged to B and E ---
$ws.Columns.Item(2).ColumnWidth = 14.3
$ws.Columns.Item(5).ColumnWidth = 33.3

# --- Move the "Process Type" list validation from E2:E16 to B2:B16 ---
$ws.Cells.Validation.Delete()
$rngB = $ws.Range("B2:B16")
$rngB.Validation.Add(3, 2, 1, "=Validate!`$B`$2:`$B`$3")
$rngB.Validation.ErrorTitle = "Process Type"
$rngB.Validation.ErrorMessage = "Invalid Process Type"
$rngB.Validation.ShowError = $true
$rngB.Validation.ShowInput = $true
$rngB.Validation.IgnoreBlank = $true

# --- Update the "cgam_processes" sheet-scoped defined name range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Processes!cgam_processes") {
        $n.RefersTo = "=Processes!`$A`$1:`$D`$7"
    }
}

# --- Update the sheet selection/active cell ---
$ws.Activate() | Out-Null
$ws.Range("E1:E16").Select() | Out-Null
